$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows (columns: A=date serial, B=nuovi pos., C=somma mobile 7gg., D=somma mobile 7gg. per 100mila abitanti)
# Update through 2021-12-08 (8/12) as per commit message.
$newData = @(
    @(44460,0,11,45.60152557831026),
    @(44461,1,12,49.7471188127021),
    @(44462,0,12,49.7471188127021),
    @(44463,1,9,37.31033910952657),
    @(44464,0,9,37.31033910952657),
    @(44465,0,6,24.87355940635105),
    @(44466,1,3,12.43677970317552),
    @(44467,3,6,24.87355940635105),
    @(44468,0,5,20.72796617195921),
    @(44469,4,9,37.31033910952657),
    @(44470,0,8,33.16474587513473),
    @(44471,1,9,37.31033910952657),
    @(44472,0,9,37.31033910952657),
    @(44473,1,9,37.31033910952657),
    @(44474,0,6,24.87355940635105),
    @(44475,1,7,29.01915264074289),
    @(44476,1,4,16.58237293756737),
    @(44477,0,4,16.58237293756737),
    @(44478,1,4,16.58237293756737),
    @(44479,4,8,33.16474587513473),
    @(44480,1,8,33.16474587513473),
    @(44481,3,11,45.60152557831026),
    @(44482,0,10,41.45593234391841),
    @(44483,1,10,41.45593234391841),
    @(44484,4,14,58.03830528148578),
    @(44485,0,13,53.89271204709394),
    @(44486,2,11,45.60152557831026),
    @(44487,1,11,45.60152557831026),
    @(44488,1,9,37.31033910952657),
    @(44489,0,9,37.31033910952657),
    @(44490,0,8,33.16474587513473),
    @(44491,3,7,29.01915264074289),
    @(44492,0,7,29.01915264074289),
    @(44493,1,6,24.87355940635105),
    @(44494,0,5,20.72796617195921),
    @(44495,0,4,16.58237293756737),
    @(44496,0,4,16.58237293756737),
    @(44497,0,4,16.58237293756737),
    @(44498,5,6,24.87355940635105),
    @(44499,0,6,24.87355940635105),
    @(44500,2,7,29.01915264074289),
    @(44501,3,10,41.45593234391841),
    @(44502,2,12,49.7471188127021),
    @(44503,0,12,49.7471188127021),
    @(44504,0,12,49.7471188127021),
    @(44505,4,11,45.60152557831026),
    @(44506,0,11,45.60152557831026),
    @(44507,5,14,58.03830528148578),
    @(44508,1,12,49.7471188127021),
    @(44509,3,13,53.89271204709394),
    @(44510,0,13,53.89271204709394),
    @(44511,3,16,66.32949175026947),
    @(44512,3,15,62.18389851587763),
    @(44513,0,15,62.18389851587763),
    @(44514,3,13,53.89271204709394),
    @(44515,0,12,49.7471188127021),
    @(44516,8,17,70.47508498466131),
    @(44517,0,17,70.47508498466131),
    @(44518,3,17,70.47508498466131),
    @(44519,8,22,91.20305115662052),
    @(44520,0,22,91.20305115662052),
    @(44521,5,24,99.4942376254042),
    @(44522,5,29,120.2222037973634),
    @(44523,6,27,111.9310173285797),
    @(44524,14,41,169.9693226100655),
    @(44525,1,39,161.6781361412818),
    @(44526,15,46,190.6972887820247),
    @(44527,4,50,207.2796617195921),
    @(44528,13,58,240.4444075947268),
    @(44529,2,55,228.0076278915513),
    @(44530,11,60,248.7355940635105),
    @(44531,1,47,194.8428820164166),
    @(44532,1,47,194.8428820164166),
    @(44533,34,66,273.6091534698615),
    @(44534,1,63,261.172373766686),
    @(44535,9,59,244.5900008291186),
    @(44536,9,66,273.6091534698615),
    @(44537,19,74,306.7738993449963),
    @(44538,0,73,302.6283061106044)
)

$startRow = 386
for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    # Copy formatting (style) from the row above so the new row matches the
    # existing table's look (bold/centered/bordered date column, datetime format).
    $ws.Range("A" + ($r - 1) + ":D" + ($r - 1)).Copy($ws.Range("A" + $r + ":D" + $r))

    $vals = $newData[$i]
    $ws.Range("A" + $r).Value = $vals[0]
    $ws.Range("B" + $r).Value = $vals[1]
    $ws.Range("C" + $r).Value = $vals[2]
    $ws.Range("D" + $r).Value = $vals[3]
}
